$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Overall Demand" row (row 2) with the sums of the per-category values
# for each material column (B, C, D).
$ws.Range("B2").Value = 202.5141002774836
$ws.Range("C2").Value = 2505.639554614019
$ws.Range("D2").Value = 1271.565699860239
